$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row above row 705; this pushes the old rows 705:746
# down to 706:747 (dimension grows from D746 to D747).
$ws.Rows.Item(705).Insert()

# Column A (date) would be auto-parsed into a date serial if we assigned
# the literal text directly, so copy an existing "2026/01/24" text cell
# (A701) into the new row to preserve the inline-string storage.
$ws.Cells.Item(701, 1).Copy() | Out-Null
$ws.Cells.Item(705, 1).PasteSpecial() | Out-Null

$ws.Cells.Item(705, 2).Value = "土"
$ws.Cells.Item(705, 3).Value = 19
$ws.Cells.Item(705, 4).Value = 201
